$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts existing rows 10-74 down to 11-75
# (matches the diff: old D10 -> new D11, old D74 -> new D75, dimension A1:R74 -> A1:R75)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with this week's data
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44819
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 100112022
$ws.Range("G10").Value = "Arveja Verde"
$ws.Range("H10").Value = "Perfection"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 25000
$ws.Range("L10").Value = 28000
$ws.Range("M10").Value = 26500
$ws.Range("N10").Value = "$/malla 25 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 1060
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
